# Revert "made changes to template files to include ID to blog page"
# This reverts commit 0a304ed20a5f255a9d1365649c3b37b0f7831aac.
#
# The prior commit had reworded four test-scenario names in column A
# (and added an "ID" concept to the blog-page test steps). This change
# reverts those four cells back to their original wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "Viewing list of posts with a particular category"
$ws.Range("A7").Value = "Leaving a comment with a name of more than 60 characters"
$ws.Range("A8").Value = "Leaving a comment without entering the name"
$ws.Range("A9").Value = "Leaving a comment without entering the comment body"

# Restore the previous selection / scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B10").Select()
